$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.572.91"
$ws.Range("E2").Value = "  -3.31%  "

$ws.Range("D3").Value = "1.782.16"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4312"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3685"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8583"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "1.773.32"
$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.490"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.281"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06955"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008754"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").Value = "26.567.74"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.134"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.33%  "

$ws.Range("D24").Value = "1.976.60"
$ws.Range("E24").Value = "  -3.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.877"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.791"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09003"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7324"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.127"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.370"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.78%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.759"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.59%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05201"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01898"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4972"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1625"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.602"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.377"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.110"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.42%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06224"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.616"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.776"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.72%  "
